$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F3").Value = "Es soll eine Ungestörte Nutzererfahrung geben."
$ws.Range("G3").Value = "Carlo Gliech"
$ws.Range("I3").Copy()
$ws.Range("J3").PasteSpecial(-4122)
$ws.Range("J3").Value = 45594
$ws.Range("K3").Value = "Abgeschlossen"
$ws.Rows.Item(3).RowHeight = 38.25
$ws.Range("Q4").Select()
